# Updates cryptos list price/volume cells (commit: "Updated cryptos list on Sun Jun  9 08:54:17 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.332.09"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "3.687.32"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "681.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.12"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.440"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000232"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "4.308.23"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.50"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.75%  "
$ws.Range("D15").Value = "3.692.58"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "69.315.61"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.08"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").Value = "  -1.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "469.18"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.95"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.655"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.85"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "3.833.61"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -5.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.92"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.91%  "
$ws.Range("E29").Value = "  -1.30%  "
$ws.Range("E30").Value = "  -4.51%  "
$ws.Range("E31").Value = "  -2.88%  "
$ws.Range("E32").Value = "  -3.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.94"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").Value = "3.675.14"
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("E36").Value = "  -5.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.28"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.28"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("E40").Value = "  -2.98%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "170.41"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.91%  "
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.63"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("E46").Value = "  -5.83%  "
$ws.Range("E47").Value = "  -2.59%  "
$ws.Range("E48").Value = "  -3.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.31"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000276"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.34%  "
$ws.Range("E51").Value = "  -3.67%  "
